$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new rows of notes/log entries to the sheet.
# Order of assignment matters: it controls the order in which new
# strings are appended to the shared string table, so it is kept in
# sync with the target shared-string index order (6 .. 23), with the
# D8 "文本框" cell assigned last (it ends up as shared string 23).
$ws.Range("A9").Value = "12061049_谢何涛_学生成绩相关性分析与系统设计实现.docx"
$ws.Range("B9").Value = "这篇修改之后没什么问题"
$ws.Range("B10").Value = "关键词部分目前还没想好"
$ws.Range("B11").Value = "关于封面的单位代码的缩进有一些问题"

$ws.Range("A14").Value = "12061053-李子靓-面向产品评论的情感要素抽取及情感倾向性分析(3).docx"
$ws.Range("B14").Value = "左侧缩进还未修改"
$ws.Range("B15").Value = "段首有Tab还未决定需要改否"

$ws.Range("A18").Value = "1145黄强_毕业论文.docx"
$ws.Range("B18").Value = "摘要字体未改正，发现rFonts标签中含有一eastAsiaTheme标签，即文章中使用主题字体"
$ws.Range("B19").Value = "部分地方含有左侧缩进和悬挂缩进"
$ws.Range("B20").Value = "tab键目前未消"
$ws.Range("B21").Value = "目录地方第三级目录标题采用了(1)检测为一级目录标题"

$ws.Range("A23").Value = "第二稿-62231365-钟华-二班-智能手机邮件收发的研究与实现(Android).docx"
$ws.Range("B24").Value = "英文摘要部分顺序不对导致定位错误"
$ws.Range("B23").Value = "部分地方有右侧缩进与左侧缩进"
$ws.Range("B25").Value = "关于列表项的编号的字体样式存在于number.xml文件中"
$ws.Range("B26").Value = "关键词还未改正"

$ws.Range("D8").Value = "文本框"

# Column A needs to widen to fit the new, longer filenames.
$ws.Columns("A").ColumnWidth = 56

# Leave the active selection on the newly added D8 cell.
$ws.Range("D8").Select() | Out-Null
